$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 7138056
$ws.Range("I33").Value = 14275902
$ws.Range("K33").Value = 14275902
$ws.Range("M33").Value = -14275673
$ws.Range("H41").Value = 368.5
$ws.Range("I41").Value = 394.53333
$ws.Range("J41").Value = 290.4
$ws.Range("K41").Value = 394.53333
$ws.Range("L41").Value = 290.4
$ws.Range("M41").Value = 45.46667000000002
$ws.Range("N41").Value = -1170.4
$ws.Range("H64").Value = 3820.05
$ws.Range("I64").Value = 3599.7778
$ws.Range("J64").Value = 4000.2727
$ws.Range("K64").Value = 3599.7778
$ws.Range("L64").Value = 4000.2727
$ws.Range("M64").Value = -3351.7778
$ws.Range("N64").Value = -4496.2727
$ws.Range("H67").Value = 3820.05
$ws.Range("I67").Value = 3599.7778
$ws.Range("J67").Value = 4000.2727
$ws.Range("K67").Value = 3599.7778
$ws.Range("L67").Value = 4000.2727
$ws.Range("M67").Value = -2741.7778
$ws.Range("N67").Value = -5716.2727
$ws.Range("H74").Value = 4266.6665
$ws.Range("I74").Value = 4120
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 4120
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -3184
$ws.Range("N74").Value = -6872
$ws.Range("H76").Value = 4647.9165
$ws.Range("I76").Value = 3829.4119
$ws.Range("J76").Value = 6635.7144
$ws.Range("K76").Value = 3829.4119
$ws.Range("L76").Value = 6635.7144
$ws.Range("M76").Value = -3514.4119
$ws.Range("N76").Value = -7265.7144
$ws.Range("H77").Value = 4266.6665
$ws.Range("I77").Value = 4120
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 20600
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -15920
$ws.Range("N77").Value = -34360
$ws.Range("H79").Value = 4647.9165
$ws.Range("I79").Value = 3829.4119
$ws.Range("J79").Value = 6635.7144
$ws.Range("K79").Value = 3829.4119
$ws.Range("L79").Value = 6635.7144
$ws.Range("M79").Value = -2737.4119
$ws.Range("N79").Value = -8819.714400000001
$ws.Range("H96").Value = 1058.4445
$ws.Range("I96").Value = 1058.4445
$ws.Range("K96").Value = 3175.3335
$ws.Range("M96").Value = -1802.3335
$ws.Range("H132").Value = 801.8
$ws.Range("I132").Value = 729.1818
$ws.Range("K132").Value = 2187.5454
$ws.Range("M132").Value = 342.4546

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 8226.857
$ws.Range("I26").Value = 717.6
$ws.Range("J26").Value = 27000
$ws.Range("K26").Value = 717.6
$ws.Range("L26").Value = 27000
$ws.Range("M26").Value = -387.6
$ws.Range("N26").Value = -27660
$ws.Range("H39").Value = 17000
$ws.Range("I39").Value = 17000
$ws.Range("K39").Value = 17000
$ws.Range("M39").Value = -16480
$ws.Range("H63").Value = 2647.3333
$ws.Range("I63").Value = 2716.8
$ws.Range("J63").Value = 2300
$ws.Range("K63").Value = 2716.8
$ws.Range("L63").Value = 2300
$ws.Range("M63").Value = -2030.8
$ws.Range("N63").Value = -3672
$ws.Range("H66").Value = 2647.3333
$ws.Range("I66").Value = 2716.8
$ws.Range("J66").Value = 2300
$ws.Range("K66").Value = 13584
$ws.Range("L66").Value = 11500
$ws.Range("M66").Value = -10152
$ws.Range("N66").Value = -18364
$ws.Range("H128").Value = 49166.668
$ws.Range("J128").Value = 49166.668
$ws.Range("L128").Value = 49166.668
$ws.Range("N128").Value = -59126.668

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2035.8667
$ws.Range("I99").Value = 1182.3334
$ws.Range("K99").Value = 1182.3334
$ws.Range("M99").Value = 315.6666

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -4713
$ws.Range("N16").ClearContents()
$ws.Range("H50").Value = 23874.5
$ws.Range("J50").Value = 23874.5
$ws.Range("L50").Value = 23874.5
$ws.Range("N50").Value = -25124.5
$ws.Range("H60").Value = 77666.5
$ws.Range("J60").Value = 62000
$ws.Range("L60").Value = 62000
$ws.Range("N60").Value = -63022
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2830
$ws.Range("N113").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4700.3076
$ws.Range("I70").Value = 2204
$ws.Range("K70").Value = 6612
$ws.Range("M70").Value = -6297
$ws.Range("H73").Value = 4700.3076
$ws.Range("I73").Value = 2204
$ws.Range("K73").Value = 6612
$ws.Range("M73").Value = -5520
$ws.Range("H87").Value = 5289.625
$ws.Range("I87").Value = 3267
$ws.Range("J87").Value = 6503.2
$ws.Range("K87").Value = 9801
$ws.Range("L87").Value = 19509.6
$ws.Range("M87").Value = -8553
$ws.Range("N87").Value = -22005.6
$ws.Range("H90").Value = 5289.625
$ws.Range("I90").Value = 3267
$ws.Range("J90").Value = 6503.2
$ws.Range("K90").Value = 29403
$ws.Range("L90").Value = 58528.8
$ws.Range("M90").Value = -23163
$ws.Range("N90").Value = -71008.79999999999
$ws.Range("H107").Value = 632.0192
$ws.Range("J107").Value = 1900
$ws.Range("L107").Value = 5700
$ws.Range("N107").Value = -9540
$ws.Range("H133").Value = 5021.407
$ws.Range("I133").Value = 3336.25
$ws.Range("J133").Value = 5730.9473
$ws.Range("K133").Value = 10008.75
$ws.Range("L133").Value = 17192.8419
$ws.Range("M133").Value = -4948.75
$ws.Range("N133").Value = -27312.8419

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5899.561
$ws.Range("I70").Value = 5660.25
$ws.Range("K70").Value = 5660.25
$ws.Range("M70").Value = -5390.25
$ws.Range("H73").Value = 5899.561
$ws.Range("I73").Value = 5660.25
$ws.Range("K73").Value = 5660.25
$ws.Range("M73").Value = -4724.25

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -3590
$ws.Range("H27").Value = 3000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3000
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -3214
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H100").Value = 4138.206
$ws.Range("I100").Value = 3426.0435
$ws.Range("K100").Value = 3426.0435
$ws.Range("M100").Value = -2885.0435

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 58500
$ws.Range("J99").Value = 58500
$ws.Range("L99").Value = 58500
$ws.Range("N99").Value = -64490
